$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 492.27118
$ws.Range("I17").Value = 426
$ws.Range("K17").Value = 1278
$ws.Range("M17").Value = -1110

$ws.Range("H19").Value = 2349.625
$ws.Range("I19").Value = 1000
$ws.Range("J19").Value = 2542.4285
$ws.Range("K19").Value = 1000
$ws.Range("L19").Value = 2542.4285
$ws.Range("M19").Value = -825
$ws.Range("N19").Value = -2892.4285

$ws.Range("H53").Value = 4317.6
$ws.Range("I53").Value = 4764.222
$ws.Range("J53").Value = 298
$ws.Range("K53").Value = 4764.222
$ws.Range("L53").Value = 298
$ws.Range("M53").Value = -4127.222
$ws.Range("N53").Value = -1572

$ws.Range("H62").Value = 52651744
$ws.Range("J62").Value = 32850
$ws.Range("L62").Value = 32850
$ws.Range("N62").Value = -34098

$ws.Range("H64").Value = 36910.617
$ws.Range("J64").Value = 7070
$ws.Range("L64").Value = 7070
$ws.Range("N64").Value = -7566

$ws.Range("H65").Value = 52651744
$ws.Range("J65").Value = 32850
$ws.Range("L65").Value = 164250
$ws.Range("N65").Value = -170490

$ws.Range("H67").Value = 36910.617
$ws.Range("J67").Value = 7070
$ws.Range("L67").Value = 7070
$ws.Range("N67").Value = -8786

$ws.Range("H74").Value = 4663.7
$ws.Range("I74").Value = 4335.125
$ws.Range("K74").Value = 4335.125
$ws.Range("M74").Value = -3399.125

$ws.Range("H77").Value = 4663.7
$ws.Range("I77").Value = 4335.125
$ws.Range("K77").Value = 21675.625
$ws.Range("M77").Value = -16995.625

$ws.Range("H100").Value = 83847.28999999999
$ws.Range("I100").Value = 54532.91
$ws.Range("J100").Value = 191333.33
$ws.Range("K100").Value = 54532.91
$ws.Range("L100").Value = 191333.33
$ws.Range("M100").Value = -53991.91
$ws.Range("N100").Value = -192415.33

$ws.Range("H112").Value = 34871.434
$ws.Range("I112").Value = 1294.4445
$ws.Range("K112").Value = 3883.3335
$ws.Range("M112").Value = -2775.3335

$ws.Range("H140").Value = 78000
$ws.Range("J140").Value = 78000
$ws.Range("L140").Value = 78000
$ws.Range("N140").Value = -88360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4790.898
$ws.Range("I32").Value = 5105.9775
$ws.Range("J32").Value = 1246.25
$ws.Range("K32").Value = 5105.9775
$ws.Range("L32").Value = 1246.25
$ws.Range("M32").Value = -4818.9775
$ws.Range("N32").Value = -1820.25

$ws.Range("H63").Value = 3168.0908
$ws.Range("I63").Value = 3233.2222
$ws.Range("K63").Value = 3233.2222
$ws.Range("M63").Value = -2547.2222

$ws.Range("H66").Value = 3168.0908
$ws.Range("I66").Value = 3233.2222
$ws.Range("K66").Value = 16166.111
$ws.Range("M66").Value = -12734.111

$ws.Range("H74").Value = 4277.643
$ws.Range("I74").Value = 2414.7896
$ws.Range("K74").Value = 2414.7896
$ws.Range("M74").Value = -1540.7896

$ws.Range("H77").Value = 4277.643
$ws.Range("I77").Value = 2414.7896
$ws.Range("K77").Value = 12073.948
$ws.Range("M77").Value = -7705.948

$ws.Range("H98").Value = 140666.67
$ws.Range("J98").Value = 140666.67
$ws.Range("L98").Value = 140666.67
$ws.Range("N98").Value = -146656.67

$ws.Range("H102").Value = 9260.485000000001
$ws.Range("I102").Value = 9626.130999999999
$ws.Range("K102").Value = 9626.130999999999
$ws.Range("M102").Value = -8004.130999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7035.2
$ws.Range("J86").Value = 4072.3333
$ws.Range("L86").Value = 4072.3333
$ws.Range("N86").Value = -6318.3333

$ws.Range("H89").Value = 7035.2
$ws.Range("J89").Value = 4072.3333
$ws.Range("L89").Value = 20361.6665
$ws.Range("N89").Value = -31593.6665

$ws.Range("H134").Value = 2586.8572
$ws.Range("I134").Value = 2149.0435
$ws.Range("K134").Value = 6447.130500000001
$ws.Range("M134").Value = -3912.130500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1036.1177
$ws.Range("I7").Value = 1595.9
$ws.Range("K7").Value = 1595.9
$ws.Range("M7").Value = -1482.9

$ws.Range("H31").Value = 2991
$ws.Range("I31").Value = 1321.3334
$ws.Range("J31").Value = 8000
$ws.Range("K31").Value = 1321.3334
$ws.Range("L31").Value = 8000
$ws.Range("M31").Value = -1026.3334
$ws.Range("N31").Value = -8590

$ws.Range("H34").Value = 2991
$ws.Range("I34").Value = 1321.3334
$ws.Range("J34").Value = 8000
$ws.Range("K34").Value = 1321.3334
$ws.Range("L34").Value = 8000
$ws.Range("M34").Value = -1119.3334
$ws.Range("N34").Value = -8404

$ws.Range("H99").Value = 11625580
$ws.Range("I99").Value = 14530287
$ws.Range("J99").Value = 6750
$ws.Range("K99").Value = 14530287
$ws.Range("L99").Value = 6750
$ws.Range("M99").Value = -14528789
$ws.Range("N99").Value = -9746

$ws.Range("H126").Value = 11625580
$ws.Range("I126").Value = 14530287
$ws.Range("J126").Value = 6750
$ws.Range("K126").Value = 43590861
$ws.Range("L126").Value = 20250
$ws.Range("M126").Value = -43588391
$ws.Range("N126").Value = -25190

$ws.Range("H132").Value = 29182.477
$ws.Range("I132").Value = 9239.5
$ws.Range("K132").Value = 27718.5
$ws.Range("M132").Value = -25188.5

$ws.Range("H139").Value = 53699.4
$ws.Range("J139").Value = 53699.4
$ws.Range("L139").Value = 53699.4
$ws.Range("N139").Value = -63979.4

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 3140.6365
$ws.Range("I137").Value = 2359.7144
$ws.Range("J137").Value = 4507.25
$ws.Range("K137").Value = 7079.1432
$ws.Range("L137").Value = 13521.75
$ws.Range("M137").Value = -1979.1432
$ws.Range("N137").Value = -23721.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5800
$ws.Range("I80").Value = 6000
$ws.Range("K80").Value = 6000
$ws.Range("M80").Value = -5002

$ws.Range("H83").Value = 5800
$ws.Range("I83").Value = 6000
$ws.Range("K83").Value = 30000
$ws.Range("M83").Value = -25008

$ws.Range("H121").Value = 60000
$ws.Range("J121").Value = 60000
$ws.Range("L121").Value = 60000
$ws.Range("N121").Value = -63494

$ws.Range("H122").Value = 22974.727
$ws.Range("I122").Value = 20627.75
$ws.Range("J122").Value = 29233.334
$ws.Range("K122").Value = 61883.25
$ws.Range("L122").Value = 87700.00199999999
$ws.Range("M122").Value = -59433.25
$ws.Range("N122").Value = -92600.00199999999

$ws.Range("H123").Value = 29750
$ws.Range("J123").Value = 29750
$ws.Range("L123").Value = 29750
$ws.Range("N123").Value = -34650

$ws.Range("H132").Value = 3261.125
$ws.Range("I132").Value = 3469.8572
$ws.Range("K132").Value = 10409.5716
$ws.Range("M132").Value = -7879.571599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3106.889
$ws.Range("I22").Value = 3410.3333
$ws.Range("J22").Value = 2500
$ws.Range("K22").Value = 3410.3333
$ws.Range("L22").Value = 2500
$ws.Range("M22").Value = -3115.3333
$ws.Range("N22").Value = -3090

$ws.Range("H25").Value = 30007
$ws.Range("I25").Value = 30007
$ws.Range("K25").Value = 30007
$ws.Range("M25").Value = -29777

$ws.Range("H27").Value = 3106.889
$ws.Range("I27").Value = 3410.3333
$ws.Range("J27").Value = 2500
$ws.Range("K27").Value = 3410.3333
$ws.Range("L27").Value = 2500
$ws.Range("M27").Value = -3303.3333
$ws.Range("N27").Value = -2714

$ws.Range("H33").Value = 10011666
$ws.Range("I33").Value = 10011666
$ws.Range("K33").Value = 10011666
$ws.Range("M33").Value = -10011376

$ws.Range("H46").Value = 3187.8823
$ws.Range("I46").Value = 839.9
$ws.Range("J46").Value = 6542.143
$ws.Range("K46").Value = 839.9
$ws.Range("L46").Value = 6542.143
$ws.Range("M46").Value = -651.9
$ws.Range("N46").Value = -6918.143

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 10571.143
$ws.Range("I29").Value = 6800
$ws.Range("K29").Value = 6800
$ws.Range("M29").Value = -6510

$ws.Range("H75").Value = 21750
$ws.Range("I75").Value = 19000
$ws.Range("K75").Value = 19000
$ws.Range("M75").Value = -18064

$ws.Range("H78").Value = 21750
$ws.Range("I78").Value = 19000
$ws.Range("K78").Value = 57000
$ws.Range("M78").Value = -52320

$ws.Range("H126").Value = 21479.682
$ws.Range("I126").Value = 25297.824
$ws.Range("J126").Value = 8498
$ws.Range("K126").Value = 75893.47200000001
$ws.Range("L126").Value = 25494
$ws.Range("M126").Value = -73423.472
$ws.Range("N126").Value = -30434
